$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4069.44
$ws.Range("I70").Value = 4656.25
$ws.Range("J70").Value = 3957.6667
$ws.Range("K70").Value = 13968.75
$ws.Range("L70").Value = 11873.0001
$ws.Range("M70").Value = -13698.75
$ws.Range("N70").Value = -12413.0001
$ws.Range("H73").Value = 4069.44
$ws.Range("I73").Value = 4656.25
$ws.Range("J73").Value = 3957.6667
$ws.Range("K73").Value = 13968.75
$ws.Range("L73").Value = 11873.0001
$ws.Range("M73").Value = -13032.75
$ws.Range("N73").Value = -13745.0001
$ws.Range("H96").Value = 656.3
$ws.Range("I96").Value = 437
$ws.Range("J96").Value = 802.5
$ws.Range("K96").Value = 1311
$ws.Range("L96").Value = 2407.5
$ws.Range("M96").Value = 62
$ws.Range("N96").Value = -5153.5
$ws.Range("H103").Value = 1470.8182
$ws.Range("I103").Value = 1495.1666
$ws.Range("J103").Value = 1441.6
$ws.Range("K103").Value = 4485.4998
$ws.Range("L103").Value = 4324.799999999999
$ws.Range("M103").Value = -3899.4998
$ws.Range("N103").Value = -5496.799999999999
$ws.Range("H112").Value = 2838.7307
$ws.Range("I112").Value = 1763.3334
$ws.Range("J112").Value = 2979
$ws.Range("K112").Value = 5290.0002
$ws.Range("L112").Value = 8937
$ws.Range("M112").Value = -4182.0002
$ws.Range("N112").Value = -11153
$ws.Range("H113").Value = 113976.664
$ws.Range("I113").Value = 2697
$ws.Range("J113").Value = 253076.25
$ws.Range("K113").Value = 2697
$ws.Range("L113").Value = 253076.25
$ws.Range("M113").Value = 557
$ws.Range("N113").Value = -259584.25
$ws.Range("H125").Value = 2257.75
$ws.Range("I125").Value = 1343.6666
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 12092.9994
$ws.Range("L125").Value = 45000
$ws.Range("M125").Value = -9632.999400000001
$ws.Range("N125").Value = -49920
$ws.Range("H132").Value = 1027.5834
$ws.Range("I132").Value = 1027.5834
$ws.Range("K132").Value = 3082.7502
$ws.Range("M132").Value = -552.7501999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 566
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888
$ws.Range("H63").Value = 8165.5713
$ws.Range("J63").Value = 9234.739
$ws.Range("L63").Value = 9234.739
$ws.Range("N63").Value = -10606.739
$ws.Range("H66").Value = 8165.5713
$ws.Range("J66").Value = 9234.739
$ws.Range("L66").Value = 46173.695
$ws.Range("N66").Value = -53037.695
$ws.Range("H74").Value = 260387.11
$ws.Range("J74").Value = 5612.375
$ws.Range("L74").Value = 5612.375
$ws.Range("N74").Value = -7360.375
$ws.Range("H77").Value = 260387.11
$ws.Range("J77").Value = 5612.375
$ws.Range("L77").Value = 28061.875
$ws.Range("N77").Value = -36797.875
$ws.Range("H110").Value = 22727928
$ws.Range("I110").Value = 27778454
$ws.Range("J110").Value = 563
$ws.Range("K110").Value = 27778454
$ws.Range("L110").Value = 563
$ws.Range("M110").Value = -27776409
$ws.Range("N110").Value = -4653
$ws.Range("H122").Value = 2038.375
$ws.Range("I122").Value = 2038.375
$ws.Range("K122").Value = 6115.125
$ws.Range("M122").Value = -3665.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 566
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -885
$ws.Range("H37").Value = 1105.1666
$ws.Range("I37").Value = 899
$ws.Range("K37").Value = 899
$ws.Range("M37").Value = -762
$ws.Range("H99").Value = 3822.3572
$ws.Range("I99").Value = 3346.3076
$ws.Range("K99").Value = 3346.3076
$ws.Range("M99").Value = -1848.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4696.278
$ws.Range("I99").Value = 4281
$ws.Range("J99").Value = 4903.9165
$ws.Range("K99").Value = 4281
$ws.Range("L99").Value = 4903.9165
$ws.Range("M99").Value = -2783
$ws.Range("N99").Value = -7899.9165
$ws.Range("H126").Value = 4696.278
$ws.Range("I126").Value = 4281
$ws.Range("J126").Value = 4903.9165
$ws.Range("K126").Value = 12843
$ws.Range("L126").Value = 14711.7495
$ws.Range("M126").Value = -10373
$ws.Range("N126").Value = -19651.7495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 48944140
$ws.Range("I4").Value = 95909410
$ws.Range("J4").Value = 5892639
$ws.Range("K4").Value = 287728230
$ws.Range("L4").Value = 17677917
$ws.Range("M4").Value = -287728118
$ws.Range("N4").Value = -17678141
$ws.Range("H23").Value = 175.33333
$ws.Range("J23").Value = 175.33333
$ws.Range("L23").Value = 525.99999
$ws.Range("N23").Value = -995.99999
$ws.Range("H33").Value = 115.6
$ws.Range("I33").Value = 109.26667
$ws.Range("J33").Value = 134.6
$ws.Range("K33").Value = 655.6000200000001
$ws.Range("L33").Value = 807.5999999999999
$ws.Range("M33").Value = -372.6000200000001
$ws.Range("N33").Value = -1373.6
$ws.Range("H35").Value = 1647
$ws.Range("I35").Value = 294
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 882
$ws.Range("L35").Value = 9000
$ws.Range("M35").Value = -594
$ws.Range("N35").Value = -9576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8103.125
$ws.Range("I43").Value = 1221
$ws.Range("J43").Value = 28749.5
$ws.Range("K43").Value = 1221
$ws.Range("L43").Value = 28749.5
$ws.Range("M43").Value = -1070
$ws.Range("N43").Value = -29051.5
$ws.Range("H46").Value = 35224
$ws.Range("J46").Value = 40094.668
$ws.Range("L46").Value = 40094.668
$ws.Range("N46").Value = -40406.668
$ws.Range("H126").Value = 2624.697
$ws.Range("J126").Value = 3628.6365
$ws.Range("L126").Value = 10885.9095
$ws.Range("N126").Value = -15825.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2525.5652
$ws.Range("I22").Value = 1733.0385
$ws.Range("J22").Value = 3555.85
$ws.Range("K22").Value = 1733.0385
$ws.Range("L22").Value = 3555.85
$ws.Range("M22").Value = -1438.0385
$ws.Range("N22").Value = -4145.85
$ws.Range("H27").Value = 2525.5652
$ws.Range("I27").Value = 1733.0385
$ws.Range("J27").Value = 3555.85
$ws.Range("K27").Value = 1733.0385
$ws.Range("L27").Value = 3555.85
$ws.Range("M27").Value = -1626.0385
$ws.Range("N27").Value = -3769.85
$ws.Range("H46").Value = 5726.1904
$ws.Range("I46").Value = 2743.25
$ws.Range("J46").Value = 7561.846
$ws.Range("K46").Value = 2743.25
$ws.Range("L46").Value = 7561.846
$ws.Range("M46").Value = -2555.25
$ws.Range("N46").Value = -7937.846
$ws.Range("H61").Value = 7622.2
$ws.Range("I61").Value = 7586.9165
$ws.Range("K61").Value = 7586.9165
$ws.Range("M61").Value = -7384.9165
$ws.Range("H68").Value = 1542.7142
$ws.Range("I68").Value = 1633.1666
$ws.Range("K68").Value = 1633.1666
$ws.Range("M68").Value = -884.1666
$ws.Range("H71").Value = 1542.7142
$ws.Range("I71").Value = 1633.1666
$ws.Range("K71").Value = 8165.833000000001
$ws.Range("M71").Value = -4421.833000000001
$ws.Range("H113").Value = 7622.2
$ws.Range("I113").Value = 7586.9165
$ws.Range("K113").Value = 7586.9165
$ws.Range("M113").Value = -5416.9165
$ws.Range("H122").Value = 41669860
$ws.Range("I122").Value = 45456664
$ws.Range("K122").Value = 136369992
$ws.Range("M122").Value = -136367542

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5250
$ws.Range("I62").Value = 5250
$ws.Range("K62").Value = 5250
$ws.Range("M62").Value = -4626
$ws.Range("H65").Value = 5250
$ws.Range("I65").Value = 5250
$ws.Range("K65").Value = 26250
$ws.Range("M65").Value = -23130
$ws.Range("H126").Value = 4441.6553
$ws.Range("I126").Value = 3992.889
$ws.Range("K126").Value = 11978.667
$ws.Range("M126").Value = -9508.667000000001
